$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells remain stored as text (matching original inline-string cells),
# since Excel would otherwise auto-convert numeric-looking strings to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "64.978.33"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "3.163.93"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "579.11"
$ws.Range("E5").Value = "  +1.05%  "
$ws.Range("D6").Value = "150.48"
$ws.Range("E6").Value = "  -0.79%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.163.37"
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D9").Value = "0.529"
$ws.Range("E9").Value = "  -0.15%  "
$ws.Range("E10").Value = "  -2.15%  "
$ws.Range("D11").Value = "6.13"
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("E12").Value = "  -0.60%  "
$ws.Range("D13").Value = "0.0000268"
$ws.Range("E13").Value = "  +3.58%  "
$ws.Range("D14").Value = "37.52"
$ws.Range("E14").Value = "  -1.57%  "
$ws.Range("D15").Value = "3.683.89"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("D16").Value = "64.964.09"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("D17").Value = "7.19"
$ws.Range("E17").Value = "  -1.07%  "
$ws.Range("D18").Value = "3.166.18"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").Value = "0.112"
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("D20").Value = "508.03"
$ws.Range("E20").Value = "  -2.55%  "
$ws.Range("D21").Value = "15.01"
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").Value = "0.721"
$ws.Range("E22").Value = "  -2.60%  "
$ws.Range("D23").Value = "15.32"
$ws.Range("E23").Value = "  +0.36%  "
$ws.Range("E24").Value = "  -1.20%  "
$ws.Range("D25").Value = "84.67"
$ws.Range("E25").Value = "  -0.94%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").Value = "9.01"
$ws.Range("E27").Value = "  +2.41%  "
$ws.Range("D28").Value = "2.93"
$ws.Range("E28").Value = "  -0.23%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").Value = "2.82"
$ws.Range("E30").Value = "  +5.34%  "
$ws.Range("D31").Value = "27.78"
$ws.Range("E31").Value = "  -1.04%  "
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("E33").Value = "  +1.21%  "
$ws.Range("D34").Value = "6.29"
$ws.Range("E34").Value = "  +2.14%  "
$ws.Range("D35").Value = "6.52"
$ws.Range("E35").Value = "  -1.15%  "
$ws.Range("D36").Value = "54.87"
$ws.Range("E36").Value = "  -1.68%  "
$ws.Range("D37").Value = "0.0894"
$ws.Range("D38").Value = "481.03"
$ws.Range("E38").Value = "  -1.70%  "
$ws.Range("D39").Value = "0.0419"
$ws.Range("E39").Value = "  -1.13%  "
$ws.Range("D40").Value = "2.97"
$ws.Range("E40").Value = "  -1.15%  "
$ws.Range("D41").Value = "8.81"
$ws.Range("E41").Value = "  +1.46%  "
$ws.Range("D42").Value = "3.013.28"
$ws.Range("E42").Value = "  -3.26%  "
$ws.Range("E43").Value = "  -3.96%  "
$ws.Range("D44").Value = "2.45"
$ws.Range("E44").Value = "  -0.84%  "
$ws.Range("D45").Value = "0.284"
$ws.Range("E45").Value = "  -4.95%  "
$ws.Range("D46").Value = "28.54"
$ws.Range("E46").Value = "  -2.82%  "
$ws.Range("D47").Value = "0.0₃0597"
$ws.Range("E47").Value = "  +2.59%  "
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("E49").Value = "  -1.62%  "
$ws.Range("D50").Value = "2.26"
$ws.Range("E50").Value = "  -2.36%  "
$ws.Range("D51").Value = "2.52"
$ws.Range("E51").Value = "  +16.62%  "
